# "Put out the fire" - update the roll count for Number Rolled = 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bump the instance count for roll value 2 (row 3) from 5 to 6.
$ws.Range("B3").Value = 6

# Recalculate all dependent formulas (percentages, sums, chart caches).
$excel.CalculateFullRebuild()

# Move the active cell selection to R26, matching where the editor left off.
$ws.Range("R26").Select()
